$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 11 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 11
}

# Zero out the production values in column B for rows 30-42
for ($r = 30; $r -le 42; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 0
}
